$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 140, shifting the existing rows (old 140-170) down to 141-171.
$ws.Rows("140:140").Insert()

# Populate the new row's data.
$ws.Range("A140").Value = "CourtOfficerNotes"
$ws.Range("B140").Value = "cares\Courts.xlsx"
$ws.Range("C140").Value = "CourtOfficerNotes"
$ws.Range("D140").Value = 1

# Give the new row the same "newly added" yellow highlight used on row 139 just
# above it, by copying only the used A:D columns (not the whole 16384-column row).
$ws.Range("A139:D139").Copy()
$ws.Range("A140:D140").PasteSpecial(-4122)  # xlPasteFormats

# Match the selection left after the edit.
$ws.Range("B140").Select()
